$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        $cell,
        [string]$value
    )
    # Force the cell to Text format first so values that look like numbers
    # (e.g. "1.003", "13.97") are stored as strings, matching the source
    # workbook's inlineStr cells instead of being auto-converted to numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    # Reset style back to the default "Normal" so we don't leave a stray
    # cell-format change behind (only the value itself should differ).
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "30.428.71"
Set-TextValue $ws.Range("E2") "  -1.17%  "
Set-TextValue $ws.Range("D3") "2.102.73"
Set-TextValue $ws.Range("E3") "  -0.34%  "
Set-TextValue $ws.Range("E4") "  +0.44%  "
Set-TextValue $ws.Range("D5") "333.17"
Set-TextValue $ws.Range("E5") "  +0.48%  "
Set-TextValue $ws.Range("D6") "1.003"
Set-TextValue $ws.Range("E6") "  +0.29%  "
Set-TextValue $ws.Range("D7") "0.5231"
Set-TextValue $ws.Range("E7") "  -1.33%  "
Set-TextValue $ws.Range("D8") "0.4584"
Set-TextValue $ws.Range("E8") "  +4.98%  "
Set-TextValue $ws.Range("D9") "53.67"
Set-TextValue $ws.Range("E9") "  +13.45%  "
Set-TextValue $ws.Range("D10") "0.08956"
Set-TextValue $ws.Range("D11") "1.174"
Set-TextValue $ws.Range("E11") "  +0.25%  "
Set-TextValue $ws.Range("D12") "24.30"
Set-TextValue $ws.Range("E12") "  -1.94%  "
Set-TextValue $ws.Range("D13") "2.101.10"
Set-TextValue $ws.Range("E13") "  +0.00%  "
Set-TextValue $ws.Range("D14") "6.757"
Set-TextValue $ws.Range("E14") "  +0.47%  "
Set-TextValue $ws.Range("D15") "7.826"
Set-TextValue $ws.Range("E15") "  +0.62%  "
Set-TextValue $ws.Range("D16") "96.29"
Set-TextValue $ws.Range("E16") "  -0.52%  "
Set-TextValue $ws.Range("E17") "  +0.33%  "
Set-TextValue $ws.Range("D18") "0.00001126"
Set-TextValue $ws.Range("E18") "  -0.45%  "
Set-TextValue $ws.Range("D19") "0.06621"
Set-TextValue $ws.Range("E19") "  -0.79%  "
Set-TextValue $ws.Range("D20") "19.21"
Set-TextValue $ws.Range("E20") "  +0.89%  "
Set-TextValue $ws.Range("D21") "1.002"
Set-TextValue $ws.Range("E21") "  +0.15%  "
Set-TextValue $ws.Range("D22") "6.277"
Set-TextValue $ws.Range("E22") "  -0.41%  "
Set-TextValue $ws.Range("D23") "30.499.83"
Set-TextValue $ws.Range("E23") "  -1.09%  "
Set-TextValue $ws.Range("D24") "12.28"
Set-TextValue $ws.Range("E24") "  +0.12%  "
Set-TextValue $ws.Range("D25") "2.356"
Set-TextValue $ws.Range("E25") "  +2.98%  "
Set-TextValue $ws.Range("D26") "2.345.52"
Set-TextValue $ws.Range("E26") "  -0.22%  "
Set-TextValue $ws.Range("D27") "22.26"
Set-TextValue $ws.Range("E27") "  -1.47%  "
Set-TextValue $ws.Range("D28") "2.549"
Set-TextValue $ws.Range("E28") "  -1.33%  "
Set-TextValue $ws.Range("E29") "  +0.79%  "
Set-TextValue $ws.Range("D30") "132.36"
Set-TextValue $ws.Range("E30") "  -0.50%  "
Set-TextValue $ws.Range("D31") "1.189"
Set-TextValue $ws.Range("E31") "  -0.55%  "
Set-TextValue $ws.Range("E32") "  -0.86%  "
Set-TextValue $ws.Range("D33") "1.676"
Set-TextValue $ws.Range("E33") "  +8.17%  "
Set-TextValue $ws.Range("D34") "6.129"
Set-TextValue $ws.Range("E34") "  -0.55%  "
Set-TextValue $ws.Range("D35") "3.928"
Set-TextValue $ws.Range("E35") "  +0.06%  "
Set-TextValue $ws.Range("D36") "10.42"
Set-TextValue $ws.Range("E36") "  +8.04%  "
Set-TextValue $ws.Range("D37") "0.02563"
Set-TextValue $ws.Range("E37") "  -1.30%  "
Set-TextValue $ws.Range("D38") "0.06791"
Set-TextValue $ws.Range("E38") "  +0.68%  "
Set-TextValue $ws.Range("D39") "5.514"
Set-TextValue $ws.Range("E39") "  -0.41%  "
Set-TextValue $ws.Range("D40") "12.69"
Set-TextValue $ws.Range("E40") "  +0.50%  "
Set-TextValue $ws.Range("D41") "0.2283"
Set-TextValue $ws.Range("E41") "  +0.28%  "
Set-TextValue $ws.Range("D42") "0.6845"
Set-TextValue $ws.Range("E42") "  +0.10%  "
Set-TextValue $ws.Range("D43") "1.246"
Set-TextValue $ws.Range("E43") "  +0.54%  "
Set-TextValue $ws.Range("D44") "2.339"
Set-TextValue $ws.Range("E44") "  +5.42%  "
Set-TextValue $ws.Range("D45") "1.002"
Set-TextValue $ws.Range("E45") "  +0.19%  "
Set-TextValue $ws.Range("B46") "EnergySwap"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D46") "13.97"
Set-TextValue $ws.Range("E46") "  +0.36%  "
Set-TextValue $ws.Range("B47") "Decentraland"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D47") "0.6355"
Set-TextValue $ws.Range("E47") "  -0.91%  "
Set-TextValue $ws.Range("D48") "3.649"
Set-TextValue $ws.Range("E48") "  +0.18%  "
Set-TextValue $ws.Range("D49") "0.00000000355"
Set-TextValue $ws.Range("E49") "  +24.22%  "
Set-TextValue $ws.Range("D50") "1.242"
Set-TextValue $ws.Range("E50") "  -1.43%  "
Set-TextValue $ws.Range("D51") "1.217"
Set-TextValue $ws.Range("E51") "  +2.15%  "
